$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (date 2019-03-05 / serial 43539 already present in A25):
#   temps = 2h30m, type = Développement,
#   description = new task about the JS page-building framework
$ws.Range("B25").Value = "2h30m"
$ws.Range("C25").Value = "Développement"
$ws.Range("D25").Value = "Création du framwork de construction de la page en Javascript (gestionnaire de pages)"

# Row 26: temps = 45m, type = Développement,
#   description = new task about the MWA page's topMenu
$ws.Range("B26").Value = "45m"
$ws.Range("C26").Value = "Développement"
$ws.Range("D26").Value = "Création du topMenu de la page MWA"

# The long description wraps to two lines, so the row grows taller.
$ws.Rows.Item(25).RowHeight = 30

# Move the active selection down to the next empty row (B27).
$ws.Range("B27").Select()
